# Restored from revision of admin on 04/23/2021 12:27:00 PM.TEST Author: admin. Type: SAVE.
# Sample Project / Main.xlsx -- Rules sheet: update the "Integer min" value
# for rule R30 (row 10) from 18 to 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C10").Value = 1
